# Detector-parameters workbook update:
#  - remove the duplicated H:M header block (was hidden cols 8-13)
#  - rename the "perc_real" header (col F) to "nspikes"
#  - give col G a custom width now that H:M are gone
#  - populate Date / Tmul / Absthresh values for several patient rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlShiftToLeft = -4159

# --- remove duplicate header cells H1:M1 (shift-left within the row only) ---
$ws.Range("H1:M1").Delete($xlShiftToLeft)

# --- rename header F1 from "perc_real" to "nspikes" ---
$ws.Range("F1").Value = "nspikes"

# --- column G: now the last real column, give it its own width ---
$ws.Range("G1").ColumnWidth = 14.86

# --- fill in Date (B), Tmul (C), Absthresh (D) for the affected rows ---
# row -> Tmul (C) value; Absthresh (D) is 50 for every affected row
$rowData = @{
    2  = 17
    5  = 17
    7  = 13
    8  = 17
    9  = 21
    10 = 17
    11 = 17
}

foreach ($r in $rowData.Keys) {
    $cVal = $rowData[$r]

    # Date column - copy the "Patient" cell's formatting (font/fill/border)
    # in this row, then overlay the date number format.
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($r, 2).Value = 44300
    $ws.Cells.Item($r, 2).NumberFormat = "m/d/yy"

    # Tmul column
    if ($r -eq 7) {
        # matches the source data's slightly different (header-style) formatting
        $ws.Cells.Item(1, 1).Copy()
        $ws.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)
    } else {
        $ws.Cells.Item($r, 1).Copy()
        $ws.Cells.Item($r, 3).PasteSpecial($xlPasteFormats)
    }
    $ws.Cells.Item($r, 3).Value = $cVal

    # Absthresh column
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($r, 4).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($r, 4).Value = 50
}
